$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Row 8 ---
# Column A holds a "YYYY.MM.DD"-style date kept as literal text in the
# source workbook (no number formats are even defined there). Assigning
# that string straight to .Value makes Excel's smart-parser convert it
# into a date serial, so lead with an apostrophe (forces text entry, just
# like typing it in the Excel UI) and then clear the resulting "quote
# prefix" formatting again so the cell keeps the sheet's default
# (unstyled) look, matching the rest of the table.
$ws.Range("A8").Value = "'2018.08.21"
$ws.Range("A8").ClearFormats()

$ws.Range("B8").Value = "16:54:35"
$ws.Range("C8").Value = "RS"
$ws.Range("D8").Value = 32
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 100
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = 250
$ws.Range("I8").Value = 0.1
$ws.Range("J8").Value = 0.9399999999999999
$ws.Range("K8").Value = "effective"
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 6.9
$ws.Range("N8").Value = 66
$ws.Range("O8").Value = 37.04999999999999
$ws.Range("P8").Value = 0.4798719937468688

# --- Row 9 ---
$ws.Range("A9").Value = "'2018.08.21"
$ws.Range("A9").ClearFormats()

$ws.Range("B9").Value = "16:57:13"
$ws.Range("C9").Value = "RS"
$ws.Range("D9").Value = 32
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 100
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = 250
$ws.Range("I9").Value = 0.1
$ws.Range("J9").Value = 0.96
$ws.Range("K9").Value = "effective"
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 9.08
$ws.Range("N9").Value = 118
$ws.Range("O9").Value = 36.4
$ws.Range("P9").Value = 0.5507834180841099
